$wb = $excel.ActiveWorkbook

# Values to update (row -> new value) in column F
$updates = @{
    2 = 17
    4 = 1468
    7 = 111
    9 = 257
}

# Apply the same update to both the "展览" and "全部类型" worksheets
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
